$d = $word.ActiveDocument

$pairs = @(
    @("548×2=1096", "700×7=4900"),
    @("815×6=4890", "583×3=1749"),
    @("983×9=8847", "528×2=1056"),
    @("930×6=5580", "997×5=4985"),
    @("307×3=921", "122×8=976"),
    @("203×8=1624", "708×3=2124"),
    @("261×8=2088", "612×2=1224"),
    @("126×7=882", "991×9=8919"),
    @("982×6=5892", "489×9=4401"),
    @("686×2=1372", "987×8=7896"),
    @("279×7=1953", "658×2=1316"),
    @("941×6=5646", "602×6=3612"),
    @("376×5=1880", "765×6=4590"),
    @("718×9=6462", "260×9=2340"),
    @("932×8=7456", "417×8=3336"),
    @("933×8=7464", "432×3=1296"),
    @("403×7=2821", "286×3=858"),
    @("578×9=5202", "114×6=684"),
    @("788×4=3152", "340×7=2380"),
    @("248×5=1240", "381×5=1905"),
    @("314×6=1884", "130×2=260"),
    @("344×9=3096", "192×8=1536"),
    @("200×4=800", "594×3=1782"),
    @("376×3=1128", "452×9=4068"),
    @("312×3=936", "411×2=822")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done replacing $($pairs.Count) cells"
